$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "Tasks" sheet selection (it currently holds the tab focus);
#    it will lose tabSelected once the new sheet is activated further below.
# ---------------------------------------------------------------------------
$tasksSheet = $wb.Worksheets.Item("Tasks")
$tasksSheet.Activate() | Out-Null
$tasksSheet.Range("F32").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Update the "FeedbackForms" sheet selection to a full-row selection
#    (row 1 selected via its row header).
# ---------------------------------------------------------------------------
$feedbackSheet = $wb.Worksheets.Item("FeedbackForms")
$feedbackSheet.Activate() | Out-Null
$feedbackSheet.Rows(1).Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add the new "Alerts" worksheet at the end of the workbook and populate
#    it with the "Create New Message" alert configuration reference data.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$alertsSheet = $wb.Worksheets.Add($null, $lastSheet)
$alertsSheet.Name = "Alerts"
$alertsSheet.Activate() | Out-Null

$alertsSheet.Range("A1").Value = "AlertSendMode"
$alertsSheet.Range("B1").Value = "Type"
$alertsSheet.Range("C1").Value = "Trigger"
$alertsSheet.Range("D1").Value = "Format"
$alertsSheet.Range("E1").Value = "Active"
$alertsSheet.Range("A1:E1").Interior.Color = 65535

$alertsSheet.Range("A2").Value = "[O]wner"
$alertsSheet.Range("A3").Value = "[S]pecific User (select below)"

$alertsSheet.Range("B2").Value = "Contacts"
$alertsSheet.Range("B3").Value = "Events"

$alertsSheet.Range("C2").Value = "Ownership Changed"
$alertsSheet.Range("C3").Value = "Note Added"

$alertsSheet.Range("D2").Value = "Email Alert"
$alertsSheet.Range("D3").Value = "Text Message Alert"

$alertsSheet.Range("E2").Value = "Yes"
$alertsSheet.Range("E3").Value = "No"

$alertsSheet.Range("E3").Select() | Out-Null
